# Insert two new data rows (179 and 180) above the current row 179,
# shifting the existing rows 179:282 down to 181:284, then populate the
# two new rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 179-180 (pushes old 179..282 to 181..284).
$ws.Range("A179:A180").EntireRow.Insert()

# New row 179
$ws.Cells.Item(179, 1).Value = 5
$ws.Cells.Item(179, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(179, 3).Value = "Maule"
$ws.Cells.Item(179, 4).Value = 44460
$ws.Cells.Item(179, 5).Value = 7
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100101
$ws.Cells.Item(179, 8).Value = "Berries"
$ws.Cells.Item(179, 9).Value = 100112025
$ws.Cells.Item(179, 10).Value = "Frutilla"
$ws.Cells.Item(179, 11).Value = "Sin especificar"
$ws.Cells.Item(179, 12).Value = "Especial"
$ws.Cells.Item(179, 13).Value = 60
$ws.Cells.Item(179, 14).Value = 20000
$ws.Cells.Item(179, 15).Value = 20000
$ws.Cells.Item(179, 16).Value = 20000
$ws.Cells.Item(179, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(179, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(179, 19).Value = 2857
$ws.Cells.Item(179, 20).Value = 7

# New row 180
$ws.Cells.Item(180, 1).Value = 5
$ws.Cells.Item(180, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(180, 3).Value = "Maule"
$ws.Cells.Item(180, 4).Value = 44460
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100101
$ws.Cells.Item(180, 8).Value = "Berries"
$ws.Cells.Item(180, 9).Value = 100112025
$ws.Cells.Item(180, 10).Value = "Frutilla"
$ws.Cells.Item(180, 11).Value = "Sin especificar"
$ws.Cells.Item(180, 12).Value = "Primera"
$ws.Cells.Item(180, 13).Value = 80
$ws.Cells.Item(180, 14).Value = 17000
$ws.Cells.Item(180, 15).Value = 17000
$ws.Cells.Item(180, 16).Value = 17000
$ws.Cells.Item(180, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(180, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(180, 19).Value = 2429
$ws.Cells.Item(180, 20).Value = 7
